# Adds three new typed columns (string/number/integer/boolean) to the
# header row of the first worksheet, mirroring the existing "abc"/"xyz"
# header cells (same bold/centered/wrap style), with matching cell
# comments and data-validation rules for the new number/integer/boolean
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (C1:F1) ------------------------------------------------
$ws.Range("C1").Value = "string"
$ws.Range("D1").Value = "number"
$ws.Range("E1").Value = "integer"
$ws.Range("F1").Value = "boolean"

# Match the existing header formatting (bold, centered, wrap text) used
# by A1/B1 -- copy the format from B1 rather than re-deriving it so the
# same style record is reused instead of minting new font/xf entries.
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)

# --- Cell comments ---------------------------------------------------------
$ws.Range("C1").AddComment("Any string")
$ws.Range("D1").AddComment("Any number")
$ws.Range("E1").AddComment("Any integer")
$ws.Range("F1").AddComment("Any boolean")

# --- Data validation rules for the new columns ------------------------------
# D: any decimal in [-1e+307, 1e+307]
$ws.Range("D2:D1048576").Validation.Add(2, 1, 1, "-1e+307", "1e+307")

# E: any whole number in [-2147483647, 2147483647]
$ws.Range("E2:E1048576").Validation.Add(1, 1, 1, "-2147483647", "2147483647")

# F: boolean, as a TRUE/FALSE dropdown list
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
